# Add "average / below-20:00-value" analysis columns (AA:AD) to the
# existing EGE olympiad worksheet.
#
#   AA1 / AB1 : column headers (shared strings)
#   AA2:AA92  : =AVERAGE(Bn:Yn)                      (daily mean temperature)
#   AB2:AB92  : =IF(AAn<Vn,1,0)                       (1 when the mean is
#                below the 20:00 reading, same comparison Z already made)
#   AC2       : =COUNT(AB2:AB92)                      (sanity-check count)
#   AD2       : explanatory note (shared string)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write AB1 before AA1 so the shared-string table gets the same index
# order as the authored workbook (0 = note about 20:00 value, 1 = "Ср
# значение").
$ws.Range("AB1").Value = "Среднее значение меньше, чем значение в 20:00 того дня"
$ws.Range("AA1").Value = "Ср значение"

for ($r = 2; $r -le 92; $r++) {
    $ws.Cells.Item($r, 27).Formula = "=AVERAGE(B" + $r + ":Y" + $r + ")"
    $ws.Cells.Item($r, 28).Formula = "=IF(AA" + $r + "<V" + $r + ",1,0)"
}

# Give the new average column the same one-decimal number format ("0.0")
# used by the other temperature columns.
$ws.Range("AA2:AA92").NumberFormat = "0.0"

$ws.Range("AC2").Formula = "=COUNT(AB2:AB92)"
$ws.Range("AD2").Value = "Функция СЧЁТ посчитала, сколько единичек выпало. А где я считал ""ИСТИНА"" такого уже не будет"

# Match the author's final view state: zoomed to 120% with AC3 selected.
$win = $ws.Application.ActiveWindow
$win.Zoom = 120
$ws.Range("AC3").Select()
